{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the \"qnut-installation\" doc update:\n//   - Rewrites the \"To install package, ...\" paragraph.\n//   - Rewrites the WordPress / Drupal 8 / Concrete5 bullet paths.\n//   - Moves the \"_GoBack\" bookmark so it sits between \"feature\" and \"s\"\n//     in the closing paragraph (text itself is unchanged).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Helper: find the (first) paragraph whose text contains `needle`.\nfunction findParagraph(needle) {\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(needle) !== -1) return p;\n  }\n  return null;\n}\n\n// 1) \"To install package, ...\" paragraph -> new instructions.\nconst installPara = findParagraph(\"To install package\");\nif (installPara) {\n  installPara.insertText(\n    \"To install package, copy all files and directories under the \\u2018peanut-files\\u2019 folder in the zip file to the Peanut root directory located under your CMS root directory.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) WordPress bullet.\nconst wpPara = findParagraph(\"WordPress:\");\nif (wpPara) {\n  wpPara.insertText(\"WordPress: wp-content/plugins/peanut\", Word.InsertLocation.replace);\n}\n\n// 3) Drupal 8 bullet.\nconst drupalPara = findParagraph(\"Drupal 8:\");\nif (drupalPara) {\n  drupalPara.insertText(\"Drupal 8:  modules/twoquakers/peanut\", Word.InsertLocation.replace);\n}\n\n// 4) Concrete5 bullet.\nconst concretePara = findParagraph(\"Concrete5:\");\nif (concretePara) {\n  concretePara.insertText(\"Concrete5: packages/knockout_view\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// 5) Move the \"_GoBack\" bookmark in the closing paragraph: it should end up\n//    right after \"feature\" (i.e. before the final \"s to your site...\").\n//    Text content of that paragraph is unchanged by this edit.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst closingParagraphs = body.paragraphs;\nclosingParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet closingPara = null;\nfor (const p of closingParagraphs.items) {\n  if (p.text.indexOf(\"get started adding\") !== -1) {\n    closingPara = p;\n    break;\n  }\n}\n\nif (closingPara) {\n  const searchResults = closingPara.getRange().search(\"feature\", { matchCase: true });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  if (searchResults.items.length > 0) {\n    const insertionPoint = searchResults.items[0].getRange(Word.RangeLocation.end);\n    insertionPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the \"qnut-installation\" doc update:\n#   - Rewrites the \"To install package, ...\" paragraph.\n#   - Rewrites the WordPress / Drupal 8 / Concrete5 bullet paths.\n#   - Moves the \"_GoBack\" bookmark so it sits between \"feature\" and \"s\"\n#     in the closing paragraph (text itself is unchanged).\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($para, [string]$newText) {\n    $r = $para.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = $newText\n}\n\nfunction Find-ParagraphContaining([string]$needle) {\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n$leftQuote = [char]0x2018\n$rightQuote = [char]0x2019\n\n# 1) \"To install package, ...\" paragraph -> new instructions.\n$installPara = Find-ParagraphContaining \"To install package\"\nif ($installPara -ne $null) {\n    $installText = \"To install package, copy all files and directories under the \" + $leftQuote + \"peanut-files\" + $rightQuote + \" folder in the zip file to the Peanut root directory located under your CMS root directory.\"\n    Set-ParagraphText $installPara $installText\n}\n\n# 2) WordPress bullet.\n$wpPara = Find-ParagraphContaining \"WordPress:\"\nif ($wpPara -ne $null) {\n    $wpText = \"WordPress: wp-content/plugins/peanut\"\n    Set-ParagraphText $wpPara $wpText\n}\n\n# 3) Drupal 8 bullet.\n$drupalPara = Find-ParagraphContaining \"Drupal 8:\"\nif ($drupalPara -ne $null) {\n    $drupalText = \"Drupal 8:  modules/twoquakers/peanut\"\n    Set-ParagraphText $drupalPara $drupalText\n}\n\n# 4) Concrete5 bullet.\n$concretePara = Find-ParagraphContaining \"Concrete5:\"\nif ($concretePara -ne $null) {\n    $concreteText = \"Concrete5: packages/knockout_view\"\n    Set-ParagraphText $concretePara $concreteText\n}\n\n# 5) Move the \"_GoBack\" bookmark in the closing paragraph: it should end up\n#    right after \"feature\" (i.e. before the final \"s to your site...\").\n#    Text content of that paragraph is unchanged by this edit.\n$closingPara = Find-ParagraphContaining \"get started adding\"\nif ($closingPara -ne $null) {\n    $searchRange = $closingPara.Range.Duplicate()\n    $found = $searchRange.Find.Execute(\"feature\", $true)\n    if ($found) {\n        $searchRange.Collapse(0) | Out-Null\n\n        if ($d.Bookmarks.Exists(\"_GoBack\")) {\n            $d.Bookmarks.Item(\"_GoBack\").Delete()\n        }\n        $d.Bookmarks.Add(\"_GoBack\", $searchRange) | Out-Null\n    }\n}\n\nWrite-Output \"qnut-installation edits applied\"\n"}
